$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") should carry the same header
# style as the existing header cells (e.g. H1). Copy H1's formatting into
# I1:J1 first so the style is reused (not duplicated), then set the text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data cells in row 2
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
